$wb = $excel.ActiveWorkbook

# --- Sheet "OFF": update row 3 (label "R") with Week 13 totals ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 208
$wsOff.Range("C3").Value = 158
$wsOff.Range("D3").Value = 54
$wsOff.Range("E3").Value = 24
$wsOff.Range("G3").Value = 6

# --- Sheet "DEF": update row 3 (label "R") with Week 13 totals ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 181
$wsDef.Range("C3").Value = 113
$wsDef.Range("D3").Value = 51
$wsDef.Range("E3").Value = 24
$wsDef.Range("F3").Value = 7
$wsDef.Range("G3").Value = 4
